$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Remove the hidden "_xlchart.v1.*" defined names left over from a
#    deleted/rebuilt chart (collect names first, then delete by name -
#    deleting while iterating the live collection skips entries).
# ------------------------------------------------------------------
$defNames = @()
foreach ($n in $wb.Names) {
    $defNames += $n.Name
}
foreach ($nm in $defNames) {
    $wb.Names.Item($nm).Delete() | Out-Null
}

# ------------------------------------------------------------------
# 2. Update the underlying data that feeds each sheet / chart with the
#    new measured timings.
# ------------------------------------------------------------------
$wsNumpy = $wb.Worksheets.Item("DASK NUMPY")
$wsNumpy.Range("B1").Value = 53.89
$wsNumpy.Range("B1").NumberFormat = "0.000"
$wsNumpy.Range("B2").Value = 101.604

$wsDistributed = $wb.Worksheets.Item("DASK Distributed Chunk Sizes ")
$wsDistributed.Range("B2").Value = 3.197
$wsDistributed.Range("B3").Value = 8.57
$wsDistributed.Range("B4").Value = 30.078
$wsDistributed.Range("B5").Value = 116.413

$wsLocal = $wb.Worksheets.Item("DASK Local Chunk Sizes")
$wsLocal.Range("B2").Value = 2.684
$wsLocal.Range("B3").Value = 11.852
$wsLocal.Range("B4").Value = 50.336
$wsLocal.Range("B5").Value = 250.167

# ------------------------------------------------------------------
# 3. Re-create the per-sheet selection / active-tab state recorded in
#    the workbook. Each worksheet keeps its own cursor position, so we
#    activate each sheet in turn to stamp its selection, finishing on
#    "DASK NUMPY" so it is the tab that is active when the file is
#    saved (tabSelected moves off "Data Types").
# ------------------------------------------------------------------
$wsDataTypes = $wb.Worksheets.Item("Data Types")

$wsDataTypes.Activate() | Out-Null
$wsDataTypes.Range("B1").Select() | Out-Null

$wsDistributed.Activate() | Out-Null
$wsDistributed.Range("M14").Select() | Out-Null

$wsLocal.Activate() | Out-Null
$wsLocal.Range("B7").Select() | Out-Null

$wsNumpy.Activate() | Out-Null
$wsNumpy.Range("F4").Select() | Out-Null
